$wb = $excel.ActiveWorkbook

# --- NewLoanInput sheet (sheet1) ---
$ws1 = $wb.Worksheets.Item("NewLoanInput")
$ws1.Activate()

# rename the existing "Chaithanya 123" text used for newloan name (do this
# first so the shared-string table keeps "chaithanyatest" at the old slot)
$ws1.Range("B2").Value = "chaithanyatest"

# insert a new row 7 (Firstrepaymenton / first-repayment-on date)
$ws1.Rows.Item(7).Insert()
$ws1.Range("A7").Value = "Firstrepaymenton"
$ws1.Range("B7").Value = 42036
$ws1.Range("B5").Select()

# --- Summary sheet (sheet2): selection only ---
$ws2 = $wb.Worksheets.Item("Summary")
$ws2.Activate()
$ws2.Range("A3").Select()

# --- Repayment Schedule sheet (sheet3): scroll + selection ---
$ws3 = $wb.Worksheets.Item("Repayment Schedule")
$ws3.Activate()
$ws3.Range("C11").Select()

# --- Transactions sheet (sheet4): selection + A2 value, keep as active sheet ---
$ws4 = $wb.Worksheets.Item("Transactions")
$ws4.Activate()
$ws4.Range("A2").Value = 197
$ws4.Range("B2").Select()
